# Update price tables - Ene 2026
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for each row/column (rows 2-11, columns B-H)
$data = @{
    2  = @(97287, 124958, 164447, 278785, 481255, 84731, 73087)
    3  = @(97711, 134341, 166004, 286459, 482745, 84953, 73375)
    4  = @(109915, 134615, 187054, 322428, 543553, 95468, 82764)
    5  = @(155149, 221833, 259835, 435420, 648344, 135341, 117721)
    6  = @(266772, 355351, 467701, 651888, 888349, 233751, 204007)
    7  = @(69766, 95524, 124853, 235984, 344292, 66319, 59686)
    8  = @(71074, 96468, 124814, 237928, 350683, 66002, 61602)
    9  = @(81574, 109199, 144105, 276056, 404920, 77165, 72279)
    10 = @(99419, 139757, 171643, 332353, 473824, 92703, 86834)
    11 = @(82692, 106211, 139779, 236967, 409067, 72023, 62124)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $i + 2  # Column B = 2
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}

# Update the selected cell / view state
$ws.Range("F18").Select()
